$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 347.5
$ws.Range("J19").Value = 265.66666
$ws.Range("L19").Value = 265.66666
$ws.Range("N19").Value = -615.66666
$ws.Range("H37").Value = 614
$ws.Range("I37").Value = 614
$ws.Range("K37").Value = 1842
$ws.Range("M37").Value = -1716
$ws.Range("H70").Value = 2803.4285
$ws.Range("I70").Value = 2820.75
$ws.Range("J70").Value = 2699.5
$ws.Range("K70").Value = 8462.25
$ws.Range("L70").Value = 8098.5
$ws.Range("M70").Value = -8192.25
$ws.Range("N70").Value = -8638.5
$ws.Range("H73").Value = 2803.4285
$ws.Range("I73").Value = 2820.75
$ws.Range("J73").Value = 2699.5
$ws.Range("K73").Value = 8462.25
$ws.Range("L73").Value = 8098.5
$ws.Range("M73").Value = -7526.25
$ws.Range("N73").Value = -9970.5
$ws.Range("H92").Value = 43478964
$ws.Range("I92").Value = 58824196
$ws.Range("J92").Value = 813
$ws.Range("K92").Value = 58824196
$ws.Range("L92").Value = 813
$ws.Range("M92").Value = -58822948
$ws.Range("N92").Value = -3309
$ws.Range("H115").Value = 179
$ws.Range("I115").Value = 179.66667
$ws.Range("K115").Value = 539.00001
$ws.Range("M115").Value = 1027.99999
$ws.Range("H132").Value = 2118.4285
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 998.8333
$ws.Range("I137").Value = 998.8
$ws.Range("K137").Value = 2996.4
$ws.Range("M137").Value = -446.3999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7562.5
$ws.Range("I32").Value = 5639
$ws.Range("K32").Value = 5639
$ws.Range("M32").Value = -5352
$ws.Range("H35").Value = 1706.8334
$ws.Range("I35").Value = 1706.8334
$ws.Range("K35").Value = 1706.8334
$ws.Range("M35").Value = -1300.8334
$ws.Range("H63").Value = 1825
$ws.Range("I63").Value = 1825
$ws.Range("K63").Value = 1825
$ws.Range("M63").Value = -1139
$ws.Range("H66").Value = 1825
$ws.Range("I66").Value = 1825
$ws.Range("K66").Value = 9125
$ws.Range("M66").Value = -5693
$ws.Range("H74").Value = 1485
$ws.Range("I74").Value = 1485
$ws.Range("K74").Value = 1485
$ws.Range("M74").Value = -611
$ws.Range("H77").Value = 1485
$ws.Range("I77").Value = 1485
$ws.Range("K77").Value = 7425
$ws.Range("M77").Value = -3057
$ws.Range("H110").Value = 62500224
$ws.Range("I110").Value = 299.66666
$ws.Range("K110").Value = 299.66666
$ws.Range("M110").Value = 1745.33334
$ws.Range("H122").Value = 5183
$ws.Range("I122").Value = 2774.5
$ws.Range("K122").Value = 8323.5
$ws.Range("M122").Value = -5873.5
$ws.Range("H132").Value = 10000
$ws.Range("I132").Value = 10000
$ws.Range("K132").Value = 30000
$ws.Range("M132").Value = -27470

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 895
$ws.Range("I64").Value = 898.5
$ws.Range("J64").Value = 888
$ws.Range("K64").Value = 898.5
$ws.Range("L64").Value = 888
$ws.Range("M64").Value = -673.5
$ws.Range("N64").Value = -1338
$ws.Range("H67").Value = 895
$ws.Range("I67").Value = 898.5
$ws.Range("J67").Value = 888
$ws.Range("K67").Value = 898.5
$ws.Range("L67").Value = 888
$ws.Range("M67").Value = -118.5
$ws.Range("N67").Value = -2448
$ws.Range("H86").Value = 1979.8
$ws.Range("I86").Value = 1751.8572
$ws.Range("K86").Value = 1751.8572
$ws.Range("M86").Value = -628.8571999999999
$ws.Range("H89").Value = 1979.8
$ws.Range("I89").Value = 1751.8572
$ws.Range("K89").Value = 8759.286
$ws.Range("M89").Value = -3143.286
$ws.Range("H106").Value = 25222.334
$ws.Range("J106").Value = 25222.334
$ws.Range("L106").Value = 25222.334
$ws.Range("N106").Value = -27746.334
$ws.Range("H107").Value = 35367.082
$ws.Range("I107").Value = 38128
$ws.Range("K107").Value = 38128
$ws.Range("M107").Value = -36208

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8622.75
$ws.Range("I31").Value = 7664.3335
$ws.Range("K31").Value = 7664.3335
$ws.Range("M31").Value = -7369.3335
$ws.Range("H34").Value = 8622.75
$ws.Range("I34").Value = 7664.3335
$ws.Range("K34").Value = 7664.3335
$ws.Range("M34").Value = -7462.3335
$ws.Range("H35").Value = 985
$ws.Range("I35").Value = 985
$ws.Range("K35").Value = 985
$ws.Range("M35").Value = -691
$ws.Range("H59").Value = 3161
$ws.Range("I59").Value = 3161
$ws.Range("K59").Value = 3161
$ws.Range("M59").Value = -2016
$ws.Range("H99").Value = 2006402.2
$ws.Range("I99").Value = 1258002.8
$ws.Range("K99").Value = 1258002.8
$ws.Range("M99").Value = -1256504.8
$ws.Range("H107").Value = 717.9
$ws.Range("I107").Value = 547.375
$ws.Range("J107").Value = 1400
$ws.Range("K107").Value = 547.375
$ws.Range("L107").Value = 1400
$ws.Range("M107").Value = 1372.625
$ws.Range("N107").Value = -5240
$ws.Range("H126").Value = 2006402.2
$ws.Range("I126").Value = 1258002.8
$ws.Range("K126").Value = 3774008.4
$ws.Range("M126").Value = -3771538.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 167.66667
$ws.Range("I12").Value = 10
$ws.Range("J12").Value = 199.2
$ws.Range("K12").Value = 30
$ws.Range("L12").Value = 597.5999999999999
$ws.Range("M12").Value = 143
$ws.Range("N12").Value = -943.5999999999999
$ws.Range("H14").Value = 2205
$ws.Range("I14").Value = 2205
$ws.Range("K14").Value = 6615
$ws.Range("M14").Value = -6442
$ws.Range("H68").Value = 2998.4
$ws.Range("I68").Value = 2997.6667
$ws.Range("J68").Value = 2999.5
$ws.Range("K68").Value = 8993.000100000001
$ws.Range("L68").Value = 8998.5
$ws.Range("M68").Value = -8182.000100000001
$ws.Range("N68").Value = -10620.5
$ws.Range("H71").Value = 2998.4
$ws.Range("I71").Value = 2997.6667
$ws.Range("J71").Value = 2999.5
$ws.Range("K71").Value = 26979.0003
$ws.Range("L71").Value = 26995.5
$ws.Range("M71").Value = -22923.0003
$ws.Range("N71").Value = -35107.5
$ws.Range("H86").Value = 175
$ws.Range("I86").Value = 150
$ws.Range("J86").Value = 200
$ws.Range("K86").Value = 450
$ws.Range("L86").Value = 600
$ws.Range("M86").Value = 736
$ws.Range("N86").Value = -2972
$ws.Range("H89").Value = 175
$ws.Range("I89").Value = 150
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 1350
$ws.Range("L89").Value = 1800
$ws.Range("M89").Value = 4578
$ws.Range("N89").Value = -13656
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 40000
$ws.Range("J15").Value = 40000
$ws.Range("L15").Value = 40000
$ws.Range("N15").Value = -40576
$ws.Range("H80").Value = 16124.75
$ws.Range("I80").Value = 20166.666
$ws.Range("J80").Value = 3999
$ws.Range("K80").Value = 20166.666
$ws.Range("L80").Value = 3999
$ws.Range("M80").Value = -19168.666
$ws.Range("N80").Value = -5995
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H83").Value = 16124.75
$ws.Range("I83").Value = 20166.666
$ws.Range("J83").Value = 3999
$ws.Range("K83").Value = 100833.33
$ws.Range("L83").Value = 19995
$ws.Range("M83").Value = -95841.33
$ws.Range("N83").Value = -29979
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H107").Value = 55557810
$ws.Range("I107").Value = 1750
$ws.Range("K107").Value = 1750
$ws.Range("M107").Value = 170
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -44920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1083.9642
$ws.Range("I22").Value = 936.1667
$ws.Range("J22").Value = 1350
$ws.Range("K22").Value = 936.1667
$ws.Range("L22").Value = 1350
$ws.Range("M22").Value = -641.1667
$ws.Range("N22").Value = -1940
$ws.Range("H27").Value = 1083.9642
$ws.Range("I27").Value = 936.1667
$ws.Range("J27").Value = 1350
$ws.Range("K27").Value = 936.1667
$ws.Range("L27").Value = 1350
$ws.Range("M27").Value = -829.1667
$ws.Range("N27").Value = -1564
$ws.Range("H127").Value = 82857.5
$ws.Range("J127").Value = 82857.5
$ws.Range("L127").Value = 82857.5
$ws.Range("N127").Value = -92777.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H51").Value = 17813.6
$ws.Range("I51").Value = 18023
$ws.Range("J51").Value = 17499.5
$ws.Range("K51").Value = 18023
$ws.Range("L51").Value = 17499.5
$ws.Range("M51").Value = -17513
$ws.Range("N51").Value = -18519.5
$ws.Range("H128").Value = 34042
$ws.Range("J128").Value = 34042
$ws.Range("L128").Value = 34042
$ws.Range("N128").Value = -44002
